# Applies the "output folder feature" update to the Em-3 tri229 -1.5sig workbook:
#  - Input sheet: records the missing measurement date for the sample, and
#    refreshes the re-run Monte-Carlo derived ratios/errors.
#  - Calc sheet: adds two new "Erfolgsrate" (success-rate) columns produced by
#    the new repeated-analysis output, and refreshes all re-run values.
#  - Results sheet: refreshes the re-run summary values.
#  - Constants sheet: updates the R30/29 ratio constant.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Input sheet
# ---------------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")

# New "Mess. Dat." (measurement date) entry for the sample - stored as literal
# text (matches source data layout), not an auto-converted date serial.
$wsInput.Range("D3").Value = "'03.01.2017"

# Refreshed Monte-Carlo re-run values
$wsInput.Range("L3").Value = 3.125603062778264
$wsInput.Range("M3").Value = 0.01283955242086306
$wsInput.Range("N3").Value = 3.54551924390971
$wsInput.Range("O3").Value = 1.069781537870688
$wsInput.Range("T3").Value = 0.5832367987382572
$wsInput.Range("U3").Value = 0.3212830336741677
$wsInput.Range("X3").Value = 0.002070408906527906
$wsInput.Range("Y3").Value = 0.2893874872767901

# Cosmetic column-width refresh (column N)
$wsInput.Columns.Item(14).ColumnWidth = 16.877604166666668

# ---------------------------------------------------------------------------
# Calc sheet
# ---------------------------------------------------------------------------
$wsCalc = $wb.Worksheets.Item("Calc")

# Two new trailing columns for the repeated-analyses success-rate output
$wsCalc.Range("BG1").Value = "Unkorr. Montefehler Erfolgsrate"
$wsCalc.Range("BH1").Value = "Korr. Montefehler Erfolgsrate"
$wsCalc.Range("BG2").Value = "(%)"
$wsCalc.Range("BH2").Value = "(%)"
$wsCalc.Range("BG3").Value = 100
$wsCalc.Range("BH3").Value = 100

# Refreshed Monte-Carlo re-run values
$wsCalc.Range("C3").Value = 3.54551924390971
$wsCalc.Range("D3").Value = 0.03792931029299849
$wsCalc.Range("E3").Value = 3.125603062778264
$wsCalc.Range("F3").Value = 0.0004013134437135165
$wsCalc.Range("G3").Value = 430.9612758989298
$wsCalc.Range("H3").Value = 0.09521967858441162
$wsCalc.Range("J3").Value = 115.8011804281942
$wsCalc.Range("K3").Value = 1.238819648857146
$wsCalc.Range("L3").Value = 1.598801382265693
$wsCalc.Range("M3").Value = 0.01710368201469975
$wsCalc.Range("O3").Value = 1.417913406792276
$wsCalc.Range("P3").Value = 0.0001820537351475395
$wsCalc.Range("Q3").Value = 1.057959448100371
$wsCalc.Range("R3").Value = 0.0001358372579303206
$wsCalc.Range("AC3").Value = 8.92447590554525
$wsCalc.Range("AD3").Value = 0.02867282692885593
$wsCalc.Range("AE3").Value = 0.4073732260934315
$wsCalc.Range("AF3").Value = 0.001308821059169303
$wsCalc.Range("AG3").Value = 386.2386535124539
$wsCalc.Range("AH3").Value = 1.813040746172784
$wsCalc.Range("AK3").Value = 0.3850556151513126
$wsCalc.Range("AL3").Value = 0.001238105850418715
$wsCalc.Range("AM3").Value = 0.3850556151513126
$wsCalc.Range("AN3").Value = 0.001238105850418715
$wsCalc.Range("AO3").Value = 32.6929
$wsCalc.Range("AP3").Value = 0.1274
$wsCalc.Range("AQ3").Value = 0.3896870574344889
$wsCalc.Range("AR3").Value = 0.0009969370275336692
$wsCalc.Range("AS3").Value = 3.414233510820528E-06
$wsCalc.Range("AV3").Value = 32.6389
$wsCalc.Range("AW3").Value = 0.1329
$wsCalc.Range("AX3").Value = 0.1318596370860057
$wsCalc.Range("AY3").Value = 0.4071828401079693
$wsCalc.Range("BB3").Value = 510.3475166712623
$wsCalc.Range("BC3").Value = 1.918655794583926
$wsCalc.Range("BD3").Value = 32580.9
$wsCalc.Range("BE3").Value = 65.92981854300285
$wsCalc.Range("BF3").Value = 0.4039953463076442

# Cosmetic column-width refresh
$wsCalc.Columns.Item(3).ColumnWidth = 16.877604166666668
$wsCalc.Columns.Item(7).ColumnWidth = 17.877604166666668
$wsCalc.Columns.Item(13).ColumnWidth = 19.877604166666668
$wsCalc.Columns.Item(18).ColumnWidth = 21.877604166666668
$wsCalc.Columns.Item(29).ColumnWidth = 16.877604166666668
$wsCalc.Columns.Item(30).ColumnWidth = 19.877604166666668
$wsCalc.Columns.Item(33).ColumnWidth = 18.877604166666668
$wsCalc.Columns.Item(38).ColumnWidth = 20.877604166666668
$wsCalc.Columns.Item(40).ColumnWidth = 20.877604166666668
$wsCalc.Columns.Item(44).ColumnWidth = 21.877604166666668
$wsCalc.Columns.Item(45).ColumnWidth = 22.877604166666668
$wsCalc.Columns.Item(49).ColumnWidth = 8.877604166666666
$wsCalc.Columns.Item(56).ColumnWidth = 18.877604166666668
$wsCalc.Columns.Item(57).ColumnWidth = 17.877604166666668
$wsCalc.Columns.Item(58).ColumnWidth = 18.877604166666668
$wsCalc.Columns.Item(59).ColumnWidth = 31.877604166666668
$wsCalc.Columns.Item(60).ColumnWidth = 29.877604166666668

# ---------------------------------------------------------------------------
# Results sheet
# ---------------------------------------------------------------------------
$wsResults = $wb.Worksheets.Item("Results")

$wsResults.Range("C3").Value = 1417.913406792276
$wsResults.Range("D3").Value = 0.1820537351475395
$wsResults.Range("G3").Value = 0.3850556151513126
$wsResults.Range("H3").Value = 0.001238105850418715
$wsResults.Range("I3").Value = 386.2386535124539
$wsResults.Range("J3").Value = 1.813040746172784
$wsResults.Range("M3").Value = 32.6929
$wsResults.Range("N3").Value = 0.1274
$wsResults.Range("O3").Value = 32.6389
$wsResults.Range("P3").Value = 0.1329
$wsResults.Range("Q3").Value = 510.3475166712623
$wsResults.Range("R3").Value = 1.918655794583926

# Cosmetic column-width refresh
$wsResults.Columns.Item(4).ColumnWidth = 19.877604166666668
$wsResults.Columns.Item(8).ColumnWidth = 20.877604166666668
$wsResults.Columns.Item(9).ColumnWidth = 18.877604166666668
$wsResults.Columns.Item(16).ColumnWidth = 7.877604166666667

# ---------------------------------------------------------------------------
# Constants sheet
# ---------------------------------------------------------------------------
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Range("B3").Value = 5E-05
